$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for eddf713a-...md is now "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-37-19 18:37:35"

# --- zh-cn sheet: same file's status + new handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-19 18:37:32"

# --- de-de sheet: same file's status + new handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-19 18:37:35"
